{"js": "// The published HTML->docx export used to append a small footer block after\n// the \"Requisitos\" list: a blank paragraph, a \"Ver no Jupiter ...\" line, and\n// a \"\u00a9 <year> ...\" copyright line. This rebuild drops that trailing footer\n// block entirely, leaving the blank paragraph that originally preceded the\n// page-break paragraph at the very end of the body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"...LOQ4084... (Requisito fraco)\") that is\n// immediately followed by the footer block we need to remove.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4084\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const toDelete = [];\n  // Right after the anchor: one blank paragraph, then the \"Ver no Jupiter\"\n  // paragraph, then the \"\u00a9 ...\" copyright paragraph.\n  const blank = items[anchorIndex + 1];\n  const jupiter = items[anchorIndex + 2];\n  const copyright = items[anchorIndex + 3];\n\n  if (blank && blank.text === \"\") toDelete.push(blank);\n  if (jupiter && jupiter.text.indexOf(\"Ver no Jupiter\") !== -1) toDelete.push(jupiter);\n  if (copyright && copyright.text.indexOf(\"Powered by Jekyll\") !== -1) toDelete.push(copyright);\n\n  // Delete from the bottom up so earlier indices stay valid.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# The published HTML->docx export used to append a small footer block after\n# the \"Requisitos\" list: a blank paragraph, a \"Ver no Jupiter ...\" line, and\n# a \"(c) <year> ...\" copyright line. This rebuild drops that trailing footer\n# block entirely, leaving the blank paragraph that originally preceded the\n# page-break paragraph at the very end of the body.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"...LOQ4084... (Requisito fraco)\") that is\n# immediately followed by the footer block we need to remove.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOQ4084*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    $blank = $d.Paragraphs.Item($anchorIndex + 1)\n    $jupiter = $d.Paragraphs.Item($anchorIndex + 2)\n    $copyright = $d.Paragraphs.Item($anchorIndex + 3)\n\n    # Delete from the bottom up so earlier indices stay valid.\n    if ($copyright.Range.Text -like \"*Powered by Jekyll*\") {\n        $copyright.Range.Delete()\n    }\n    if ($jupiter.Range.Text -like \"*Ver no Jupiter*\") {\n        $jupiter.Range.Delete()\n    }\n    if ($blank.Range.Text -eq \"`r\") {\n        $blank.Range.Delete()\n    }\n}\n"}
